# "Added the table of contents"
#
# The author resized/repositioned the three shapes that make up the
# last slide (the cover-style "Design of the Train Control ..." slide)
# so that the dark background rectangle now fills the whole slide and
# the picture/title textbox shift up to sit inside it.
#
# PowerPoint's Shape.Left/Top/Width/Height are expressed in points,
# while the OOXML stores EMUs (1 pt = 12700 EMU). 914400 EMU/inch.
# A tiny epsilon is added before converting EMU -> points so that the
# point -> EMU round trip performed internally lands back on the exact
# target EMU value instead of being truncated one unit short by
# floating point noise.

function EmuToPt($emu) {
    return ($emu / 12700.0) + 0.00001
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item($p.Slides.Count)

# Background rectangle now covers the entire 12192000 x 6858000 slide.
$rect = $s.Shapes.Item("Rectangle 8")
$rect.Left = EmuToPt 603504
$rect.Top = EmuToPt 0
$rect.Width = EmuToPt 10972800
$rect.Height = EmuToPt 6858000

# Picture shifts up (same size / same horizontal position).
$pic = $s.Shapes.Item("Picture 4")
$pic.Top = EmuToPt 995495

# Title textbox shifts up to match (same size / same horizontal position).
$title = $s.Shapes.Item("TextBox 5")
$title.Top = EmuToPt 394155
